$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 46002
$ws.Cells.Item(2, 2).Value = 11355.1020769428
$ws.Cells.Item(2, 3).Value = 10508.6445966388
$ws.Cells.Item(2, 4).Value = 17232.26
$ws.Cells.Item(2, 5).Value = 6834.07602783973
$ws.Cells.Item(2, 6).Value = 4.60252601993867

# Row 3
$ws.Cells.Item(3, 1).Value = 46003
$ws.Cells.Item(3, 2).Value = 11287.0042176068
$ws.Cells.Item(3, 3).Value = 9726.51135775188
$ws.Cells.Item(3, 4).Value = 11232.26
$ws.Cells.Item(3, 5).Value = 7043.73354076636
$ws.Cells.Item(3, 6).Value = 230.749370771593

# Row 4
$ws.Cells.Item(4, 1).Value = 46004
$ws.Cells.Item(4, 2).Value = 4252.87008618029
$ws.Cells.Item(4, 3).Value = 6650.24599293135
$ws.Cells.Item(4, 4).Value = 11232.26
$ws.Cells.Item(4, 5).Value = 6647.73111586245
$ws.Cells.Item(4, 6).Value = 86.0715461997416

# Row 5
$ws.Cells.Item(5, 1).Value = 46005
$ws.Cells.Item(5, 2).Value = 4174.28410869778
$ws.Cells.Item(5, 3).Value = 6833.03218006017
$ws.Cells.Item(5, 4).Value = 11232.26
$ws.Cells.Item(5, 5).Value = 6742.22186205447
$ws.Cells.Item(5, 6).Value = 97.6247517547766

# Row 6
$ws.Cells.Item(6, 1).Value = 46006
$ws.Cells.Item(6, 2).Value = 11051.1309582548
$ws.Cells.Item(6, 3).Value = 10522.7301906936
$ws.Cells.Item(6, 4).Value = 11232.26
$ws.Cells.Item(6, 5).Value = 7314.89847311247
$ws.Cells.Item(6, 6).Value = 275.223694325252

# Row 7
$ws.Cells.Item(7, 1).Value = 46007
$ws.Cells.Item(7, 2).Value = 9746.03784374469
$ws.Cells.Item(7, 3).Value = 9692.75654673297
$ws.Cells.Item(7, 4).Value = 11232.26
$ws.Cells.Item(7, 5).Value = 8055.89644031477
$ws.Cells.Item(7, 6).Value = 271.516374460322

# Row 8
$ws.Cells.Item(8, 1).Value = 46008
$ws.Cells.Item(8, 2).Value = 9746.03784374469
$ws.Cells.Item(8, 3).Value = 9176.24380258412
$ws.Cells.Item(8, 4).Value = 11232.26
$ws.Cells.Item(8, 5).Value = 8055.89644031477
$ws.Cells.Item(8, 6).Value = 249.995010120787

# Row 9
$ws.Cells.Item(9, 1).Value = 46009
$ws.Cells.Item(9, 2).Value = 9746.03784374469
$ws.Cells.Item(9, 3).Value = 9244.2086529491
$ws.Cells.Item(9, 4).Value = 11232.26
$ws.Cells.Item(9, 5).Value = 8055.89644031477
$ws.Cells.Item(9, 6).Value = 252.826878885995

# Row 10
$ws.Cells.Item(10, 1).Value = 46010
$ws.Cells.Item(10, 2).Value = 9746.03784374469
$ws.Cells.Item(10, 3).Value = 8490.44327276352
$ws.Cells.Item(10, 4).Value = 11232.26
$ws.Cells.Item(10, 5).Value = 8055.89644031477
$ws.Cells.Item(10, 6).Value = 221.419988044929

# Row 11
$ws.Cells.Item(11, 1).Value = 46011
$ws.Cells.Item(11, 2).Value = 8560.57524882407
$ws.Cells.Item(11, 3).Value = 8338.0873348358
$ws.Cells.Item(11, 4).Value = 11232.26
$ws.Cells.Item(11, 5).Value = 7663.95939515638
$ws.Cells.Item(11, 6).Value = 198.741113749674

# Row 12
$ws.Cells.Item(12, 1).Value = 46012
$ws.Cells.Item(12, 2).Value = 8469.04674334209
$ws.Cells.Item(12, 3).Value = 8444.13539390051
$ws.Cells.Item(12, 4).Value = 11232.26
$ws.Cells.Item(12, 5).Value = 7655.56137905024
$ws.Cells.Item(12, 6).Value = 202.809865539615

# Row 13
$ws.Cells.Item(13, 1).Value = 46013
$ws.Cells.Item(13, 2).Value = 9566.18303866457
$ws.Cells.Item(13, 3).Value = 9381.85017403134
$ws.Cells.Item(13, 4).Value = 11232.26
$ws.Cells.Item(13, 5).Value = 8126.77104635949
$ws.Cells.Item(13, 6).Value = 261.515050849618

# Row 14
$ws.Cells.Item(14, 1).Value = 46014
$ws.Cells.Item(14, 2).Value = 9566.18303866457
$ws.Cells.Item(14, 3).Value = 9825.97441021394
$ws.Cells.Item(14, 4).Value = 11232.26
$ws.Cells.Item(14, 5).Value = 8126.77104635949
$ws.Cells.Item(14, 6).Value = 280.020227357227

# Row 15
$ws.Cells.Item(15, 1).Value = 46015
$ws.Cells.Item(15, 2).Value = 9566.18303866457
$ws.Cells.Item(15, 3).Value = 9840.41234807663
$ws.Cells.Item(15, 4).Value = 11232.26
$ws.Cells.Item(15, 5).Value = 8126.77104635949
$ws.Cells.Item(15, 6).Value = 280.621808101505
